$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Remove the two trailing blank rows (10 and 9) - the table now
#    only spans down to row 8.
# ------------------------------------------------------------------
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(9).Delete()

# ------------------------------------------------------------------
# 2) Clear out the old client rows (2:8) so we can re-populate them
#    with the new client list.
# ------------------------------------------------------------------
$ws.Range("A2:I8").ClearContents()

# ------------------------------------------------------------------
# 3) Widen column D (Contraseña) and drop its "best fit" sizing -
#    it is now a fixed custom width.
# ------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 16.33

# ------------------------------------------------------------------
# 4) New client data (Cliente, CUIT para ingresar, CUIT representado,
#    Contraseña) - columns A:D, rows 2:8.
# ------------------------------------------------------------------
$ws.Range("A2").Value = "Alfredo Quintana"
$ws.Range("B2").Value = 20111155500
$ws.Range("C2").Value = 20111155500
$ws.Range("D2").Value = "((#Urquiza#7411))`n"

$ws.Range("A3").Value = "Burgi Omar"
$ws.Range("B3").Value = 20115710037
$ws.Range("C3").Value = 20115710037
$ws.Range("D3").Value = "Omarbu2024"

$ws.Range("A4").Value = "Fassi Alberto"
$ws.Range("B4").Value = 20170895658
$ws.Range("C4").Value = 20170895658
$ws.Range("D4").Value = "Alberto2024"

$ws.Range("A5").Value = "Florentino Rivarossa"
$ws.Range("B5").Value = 20113062518
$ws.Range("C5").Value = 20113062518
$ws.Range("D5").Value = "Estudio2024"

$ws.Range("A6").Value = "Seffino Marcelo"
$ws.Range("B6").Value = 23163038919
$ws.Range("C6").Value = 23163038919
$ws.Range("D6").Value = "Marcelo2023"

$ws.Range("A8").Value = "Gabriela Evangelina Lisi"
$ws.Range("B8").Value = 27160314066
$ws.Range("C8").Value = 27160314066
$ws.Range("D8").Value = "Gabriela2023"

$ws.Range("A7").Value = "Gabriela Evangelina Lisi y Seffino Marcelo"
$ws.Range("B7").Value = 23163038919
$ws.Range("C7").Value = 30715364170
$ws.Range("D7").Value = "Marcelo2023"

# The new password for row 2 contains an embedded line break; undo the
# automatic row-height bump that causes so the row keeps its default
# height (matches the other, untouched rows).
$ws.Rows.Item(2).EntireRow.AutoFit()

# ------------------------------------------------------------------
# 5) Formulas - columns E:I, rows 2:8 (same pattern used by the
#    existing rows, just extended down through row 8).
# ------------------------------------------------------------------
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 5).Formula = '=LEFT(CELL("filename"),FIND("[",CELL("filename"))-1)'
    $ws.Cells.Item($r, 6).Formula = '=E' + $r + '&"Deudas\"'
}

$ws.Range("G2").Formula = "=IF(B2=B1,1,0)"
$ws.Range("H2").Formula = "=IF(B2=B3,1,0)"
$ws.Range("I2").Formula = "=G2+H2"

for ($r = 3; $r -le 8; $r++) {
    $prev = $r - 1
    $next = $r + 1
    $ws.Cells.Item($r, 7).Formula = '=IF(B' + $r + '=B' + $prev + ',1,0)'
    $ws.Cells.Item($r, 8).Formula = '=IF(B' + $r + '=B' + $next + ',1,0)'
    $ws.Cells.Item($r, 9).Formula = '=G' + $r + '+H' + $r
}

# Last row (8) has no row below it, so H8 compares against the
# (now empty) row 9 - same as the original sheet's pattern.

# ------------------------------------------------------------------
# 6) Selection - the author's cursor ended up on E12.
# ------------------------------------------------------------------
$ws.Range("E12").Select()

Write-Output "done"
